$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 38 matching style/content pattern of the last existing row (37)
$ws.Range("A38").Value = "keywords"
$ws.Range("B38").Value = "/Data File Subjects/Keyword"

# Copy the style from row 37 (last data row) onto the new row 38
$ws.Range("A37:B37").Copy()
$ws.Range("A38:B38").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move selection to match the post-edit cursor location recorded in the diff
$ws.Range("B44").Select()
